$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.15 = 29015.78 pesos`n✅ 29015.78 pesos = 7.13 = 912.83 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the tasa values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 139.958
$ws2.Range("O10").Value = 4060.99
$ws2.Range("N12").Value = 4069
$ws2.Range("O12").Value = 128.01
